# "fix 131_1_20 to 131_0_20"
#
# The sheet originally had a single "index" column (C) holding the
# (5-way) MOS ordering string, e.g. "5 2 3 1 4". This edit inserts a new
# column C in front of it that records which of the 5 systems (A-E) is
# the current row's "anchor" system, shifting the old C:E columns right
# to D:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; B:E -> B, D:F (B keeps its formatting,
# and the new C inherits it, matching the original column B/C width+style
# pairing).
$ws.Columns("C").Insert()

# New column C values for rows 2-21 (one letter per row), in row order so
# the resulting shared-string table is built in the same first-seen order
# as the source edit (E, B, D, C, A).
$letters = @("E","E","B","B","D","B","C","C","C","E","A","D","A","A","B","A","B","B","B","B")

for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $letters[$i]
}
